$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 16799.166
$ws.Range("I45").Value = 2843.75
$ws.Range("K45").Value = 8531.25
$ws.Range("M45").Value = -8339.25
$ws.Range("H62").Value = 4994
$ws.Range("I62").Value = 4994
$ws.Range("K62").Value = 4994
$ws.Range("M62").Value = -4370
$ws.Range("H65").Value = 4994
$ws.Range("I65").Value = 4994
$ws.Range("K65").Value = 24970
$ws.Range("M65").Value = -21850
$ws.Range("H134").Value = 99987.5
$ws.Range("J134").Value = 99987.5
$ws.Range("L134").Value = 99987.5
$ws.Range("N134").Value = -110127.5
$ws.Range("H137").Value = 1729519.6
$ws.Range("I137").Value = 3225.422
$ws.Range("K137").Value = 9676.266
$ws.Range("M137").Value = -7126.266
$ws.Range("H138").Value = 8701.608
$ws.Range("J138").Value = 3136
$ws.Range("L138").Value = 9408
$ws.Range("N138").Value = -19688

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1738.3334
$ws.Range("I21").Value = 1738.3334
$ws.Range("K21").Value = 1738.3334
$ws.Range("M21").Value = -1364.3334
$ws.Range("H61").Value = 955382.7
$ws.Range("I61").Value = 27304.455
$ws.Range("K61").Value = 27304.455
$ws.Range("M61").Value = -27092.455
$ws.Range("H97").Value = 14897.333
$ws.Range("I97").Value = 16298.143
$ws.Range("K97").Value = 16298.143
$ws.Range("M97").Value = -15802.143
$ws.Range("H102").Value = 1188.0769
$ws.Range("I102").Value = 1203.75
$ws.Range("K102").Value = 1203.75
$ws.Range("M102").Value = 418.25
$ws.Range("H132").Value = 1550.614
$ws.Range("I132").Value = 972.8958
$ws.Range("K132").Value = 2918.6874
$ws.Range("M132").Value = -388.6873999999998
$ws.Range("H135").Value = 98855.5
$ws.Range("J135").Value = 98855.5
$ws.Range("L135").Value = 98855.5
$ws.Range("N135").Value = -108995.5
$ws.Range("H136").Value = 955382.7
$ws.Range("I136").Value = 27304.455
$ws.Range("K136").Value = 81913.365
$ws.Range("M136").Value = -79363.365

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 2000.5
$ws.Range("I12").Value = 2000.5
$ws.Range("K12").Value = 2000.5
$ws.Range("M12").Value = -1832.5
$ws.Range("H29").Value = 1694.3334
$ws.Range("I29").Value = 1694.3334
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1694.3334
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = -1405.3334
$ws.Range("H80").Value = 52633810
$ws.Range("J80").Value = 2535.3333
$ws.Range("L80").Value = 2535.3333
$ws.Range("N80").Value = -4531.3333
$ws.Range("H83").Value = 52633810
$ws.Range("J83").Value = 2535.3333
$ws.Range("L83").Value = 12676.6665
$ws.Range("N83").Value = -22660.6665
$ws.Range("H86").Value = 3670.0833
$ws.Range("I86").Value = 2255.4375
$ws.Range("J86").Value = 6499.375
$ws.Range("K86").Value = 2255.4375
$ws.Range("L86").Value = 6499.375
$ws.Range("M86").Value = -1132.4375
$ws.Range("N86").Value = -8745.375
$ws.Range("H89").Value = 3670.0833
$ws.Range("I89").Value = 2255.4375
$ws.Range("J89").Value = 6499.375
$ws.Range("K89").Value = 11277.1875
$ws.Range("L89").Value = 32496.875
$ws.Range("M89").Value = -5661.1875
$ws.Range("N89").Value = -43728.875
$ws.Range("H105").Value = 7322.1816
$ws.Range("I105").Value = 6401.7393
$ws.Range("J105").Value = 9439.2
$ws.Range("K105").Value = 6401.7393
$ws.Range("L105").Value = 9439.2
$ws.Range("M105").Value = -4654.7393
$ws.Range("N105").Value = -12933.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = $null
$ws.Range("H58").Value = 3140.3333
$ws.Range("I58").Value = 1421
$ws.Range("K58").Value = 1421
$ws.Range("M58").Value = -1218
$ws.Range("H136").Value = 3140.3333
$ws.Range("I136").Value = 1421
$ws.Range("K136").Value = 4263
$ws.Range("M136").Value = -1713

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4450.9287
$ws.Range("I32").Value = 3486.2222
$ws.Range("J32").Value = 6187.4
$ws.Range("K32").Value = 10458.6666
$ws.Range("L32").Value = 18562.2
$ws.Range("M32").Value = -10175.6666
$ws.Range("N32").Value = -19128.2
$ws.Range("H121").Value = 2398.4
$ws.Range("I121").Value = 498
$ws.Range("K121").Value = 1494
$ws.Range("M121").Value = -184
$ws.Range("H128").Value = 499999.5
$ws.Range("I128").Value = 499999.5
$ws.Range("K128").Value = 1499998.5
$ws.Range("M128").Value = -1495018.5
$ws.Range("H131").Value = 3389.9333
$ws.Range("I131").Value = 3036.3333
$ws.Range("K131").Value = 9108.999899999999
$ws.Range("M131").Value = -4068.999899999999
$ws.Range("H138").Value = 3936.611
$ws.Range("I138").Value = 4019.3125
$ws.Range("J138").Value = 3275
$ws.Range("K138").Value = 12057.9375
$ws.Range("L138").Value = 9825
$ws.Range("M138").Value = -6917.9375
$ws.Range("N138").Value = -20105
$ws.Range("H139").Value = 3311.6316
$ws.Range("I139").Value = 2193.2666
$ws.Range("K139").Value = 6579.7998
$ws.Range("M139").Value = -1439.7998
$ws.Range("H140").Value = 1609.303
$ws.Range("I140").Value = 1358.862
$ws.Range("K140").Value = 4076.586
$ws.Range("M140").Value = 1103.414

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 8747
$ws.Range("J25").Value = 8747
$ws.Range("L25").Value = 8747
$ws.Range("N25").Value = -9805
$ws.Range("H48").Value = 27750.75
$ws.Range("J48").Value = 27750.75
$ws.Range("L48").Value = 27750.75
$ws.Range("N48").Value = -28720.75
$ws.Range("H62").Value = 45998.25
$ws.Range("J62").Value = 49999.5
$ws.Range("L62").Value = 49999.5
$ws.Range("N62").Value = -51371.5
$ws.Range("H65").Value = 45998.25
$ws.Range("J65").Value = 49999.5
$ws.Range("L65").Value = 149998.5
$ws.Range("N65").Value = -156862.5
$ws.Range("H126").Value = 1862
$ws.Range("I126").Value = 1434.8334
$ws.Range("K126").Value = 4304.5002
$ws.Range("M126").Value = -1834.5002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 80159.766
$ws.Range("I13").Value = 12500
$ws.Range("J13").Value = 92461.55
$ws.Range("K13").Value = 12500
$ws.Range("L13").Value = 92461.55
$ws.Range("M13").Value = -12360
$ws.Range("N13").Value = -92741.55
$ws.Range("H43").Value = 2234230.8
$ws.Range("I43").Value = 380000
$ws.Range("K43").Value = 380000
$ws.Range("M43").Value = -379807
$ws.Range("H46").Value = 15617.333
$ws.Range("I46").Value = 16319.5
$ws.Range("K46").Value = 16319.5
$ws.Range("M46").Value = -16131.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 14995
$ws.Range("J80").Value = 14995
$ws.Range("L80").Value = 14995
$ws.Range("N80").Value = -16991
$ws.Range("H83").Value = 14995
$ws.Range("J83").Value = 14995
$ws.Range("L83").Value = 44985
$ws.Range("N83").Value = -54969
$ws.Range("H100").Value = 1722.0769
$ws.Range("I100").Value = 712
$ws.Range("K100").Value = 1424
$ws.Range("M100").Value = -883
